$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Inflammatory-Mac" (the Sending cluster label used in A4) is renamed to "MuSCs".
# Column D (Target cluster) already held the text "MuSCs" for every data row, so
# after the rename both labels collapse onto the same shared text value.
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D2:D5").Value = "MuSCs"

# Row 2 - refreshed TPM-derived statistics
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2955753333333333
$ws.Range("H2").Value = 0.8867259999999999
$ws.Range("I2").Value = 0.239018529794766
$ws.Range("J2").Value = 0.2584571780171812
$ws.Range("M2").Value = 0.1182145
$ws.Range("N2").Value = 0.236429
$ws.Range("Q2").Value = 0.03494129024233333
$ws.Range("R2").Value = 0.209647741454
$ws.Range("S2").Value = 0.239018529794766
$ws.Range("T2").Value = 0.2584571780171812

# Row 3
$ws.Range("I3").Value = 0.2089742936599006
$ws.Range("J3").Value = 0.2259695357671569
$ws.Range("M3").Value = 0.1182145
$ws.Range("N3").Value = 0.236429
$ws.Range("Q3").Value = 0.030549227519
$ws.Range("R3").Value = 0.183295365114
$ws.Range("S3").Value = 0.2089742936599006
$ws.Range("T3").Value = 0.2259695357671569

# Row 4
$ws.Range("F4").Value = 0.5
$ws.Range("G4").Value = 0.27902
$ws.Range("H4").Value = 0.55804
$ws.Range("I4").Value = 0.2256309734348681
$ws.Range("J4").Value = 0.1626539016795581
$ws.Range("M4").Value = 0.1182145
$ws.Range("N4").Value = 0.236429
$ws.Range("Q4").Value = 0.03298420979
$ws.Range("R4").Value = 0.13193683916
$ws.Range("S4").Value = 0.2256309734348681
$ws.Range("T4").Value = 0.1626539016795581

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.4036036666666667
$ws.Range("H5").Value = 1.210811
$ws.Range("I5").Value = 0.3263762031104653
$ws.Range("J5").Value = 0.3529193845361038
$ws.Range("M5").Value = 0.1182145
$ws.Range("N5").Value = 0.236429
$ws.Range("Q5").Value = 0.04771180565316667
$ws.Range("R5").Value = 0.286270833919
$ws.Range("S5").Value = 0.3263762031104653
$ws.Range("T5").Value = 0.3529193845361038
